$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.814.02"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.893.36"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "0.7572"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "240.14"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.3053"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").Value = "25.05"
$ws.Range("E9").Value = "  -7.24%  "
$ws.Range("E10").Value = "  -1.82%  "
$ws.Range("D11").Value = "0.07966"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "0.7479"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").Value = "1.898.54"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "5.198"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "91.17"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "29.826.98"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "6.034"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").Value = "13.82"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").Value = "0.000007662"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").Value = "232.81"
$ws.Range("E20").Value = "  -5.17%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "2.155.14"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").Value = "9.222"
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("D26").Value = "164.56"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").Value = "18.64"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("D28").Value = "0.1282"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "2.034"
$ws.Range("E29").Value = "  -4.50%  "
$ws.Range("D30").Value = "1.338"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").Value = "1.522"
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").Value = "4.277"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "3.998"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "0.05312"
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "1.237"
$ws.Range("E35").Value = "  -4.34%  "
$ws.Range("D36").Value = "0.7290"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "2.715"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("D38").Value = "0.01919"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "2.762"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "6.209"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "0.4409"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "72.36"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("D43").Value = "1.908"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "0.8264"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "7.588"
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").Value = "9.816"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").Value = "2.056.37"
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "35.92"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "0.05945"
$ws.Range("E51").Value = "  -0.48%  "
